# Commit: "adding new files from GD"
#
# The canonical OOXML diff shows that every paragraph in the document
# loses its explicit
#
#     <w:contextualSpacing w:val="0"/>
#
# entry from <w:pPr> -- nothing else in any paragraph's properties
# changes (borders, shading, run fonts/colors, text, etc. are all left
# untouched). In the Word object model this direct-formatting flag is
# ParagraphFormat.ContextualSpacing ("Don't add space between
# paragraphs of the same style"), so clear it on every paragraph in the
# document.

$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $para.Range.ParagraphFormat.ContextualSpacing = $false
}
